# Update "Generate Report for Handback" timestamps in the handback-status workbook.
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for b31a489c row (G4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-16 22:45:39"

# zh-cn sheet: Handoff/Handback datetimes for b31a489c row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-16 22:45:34"
$wsZhCn.Range("K4").Value = "2016-08-16 22:45:52"

# de-de sheet: same "Correspond Handoff Datetime" text as Overview's G4 (shared
# string), plus its own "Correspond Handback DateTime" (K4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-16 22:45:39"
$wsDeDe.Range("K4").Value = "2016-08-16 22:45:59"
